# Insert a new sale-item row ("NO DEPRINE 50MG 30 TAB") just above the
# "OPTIDEX -T EYE DROPS 5 ML" row (worksheet row 19), pushing every
# following row down by one. Re-number the item counter in column A,
# re-create the merged cell regions for the new row and refresh the
# totals row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Push row 19 (and everything below it) down by one row.
$ws.Rows(19).Insert()

# 2. Fill in the new row 19 with the new item's data.
$ws.Cells.Item(19, 1).Value = 13
$ws.Cells.Item(19, 3).Value = "NO DEPRINE 50MG 30 TAB"
$ws.Cells.Item(19, 8).Value = "0:0"
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 14).Value = "75.00"
$ws.Cells.Item(19, 16).Value = "75.0000"
$ws.Cells.Item(19, 17).Value = "1:0"

# 3. Re-create the merged regions for the new row 19 (mirrors every
#    other item row: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# 4. Copy formatting from the row right below (the old row 19, now
#    shifted to row 20) so the new row looks identical to its peers.
$ws.Rows(20).Copy()
$ws.Rows(19).PasteSpecial(-4122)
$ws.Cells.Item(19, 1).Value = 13
$ws.Cells.Item(19, 3).Value = "NO DEPRINE 50MG 30 TAB"
$ws.Cells.Item(19, 8).Value = "0:0"
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 14).Value = "75.00"
$ws.Cells.Item(19, 16).Value = "75.0000"
$ws.Cells.Item(19, 17).Value = "1:0"
$ws.Rows(19).RowHeight = $ws.Rows(20).RowHeight

# 5. Renumber the running item counter (column A) for every item that
#    used to sit at row 19 and below - they all moved down by one row
#    and must show one higher a sequence number than before.
for ($r = 20; $r -le 36; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value + 1
}

# 6. Update the grand-total cell (now on row 37) to include the new
#    item's selling price.
$ws.Cells.Item(37, 16).Value = $ws.Cells.Item(37, 16).Value + 75
